# SCD0172 - Melakukan Proses Pemantauan pada Menu Pipeline
#
# The workbook's data rows (row 2 / row 3) describe a test scenario around
# "Product Holding" report searching. This edit:
#   - Removes the 3rd bullet of the expected-result text (column E) and makes
#     both rows 2 and 3 share the same (shortened) expected-result text.
#   - Rewrites the test-steps text (column D) to a lower-cased wording that
#     also mentions "dengan field bulan", and makes both rows 2 and 3 share
#     this same steps text.
#   - Adjusts row 2's height to match row 3's (since its text got shorter).
#   - Moves the active selection from Q2 to E2 (and scrolls the view back
#     toward column D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCD0179")

$newSteps = "1. login digisales portal dengan sales 39798`n2. buka menu product holding`n3. isi field npp sales sendiri dengan 39798, dengan field bulan februari 2022 (karena hanya ada data ini)`n4. klik generate`n5. data muncul"
$newExpected = " - Field Nama dan Tipe Sales akan terisi secara otomatis oleh system berdasarkan field NPP`n - Field NPP, Tahun Data, Bulan Data bersifat mandatori"

# Row 2 (previously: steps without date filled in; 3-bullet expected result)
$ws.Range("D2").Value = $newSteps
$ws.Range("E2").Value = $newExpected

# Row 3 (previously: steps with bulan februari 2022; "Sesuai dengan BSDD 3.3")
$ws.Range("D3").Value = $newSteps
$ws.Range("E3").Value = $newExpected

# Row 2's wrapped text is now shorter (matches row 3), so its height shrinks
# from 89.25 to 76.5 points.
$ws.Rows.Item(2).RowHeight = 76.5

# Move the view/selection from Q2 back to E2.
$ws.Activate() | Out-Null
$ws.Range("E2").Select() | Out-Null
